# FFmpeg.xlsx commit-tracker sheet: replace the last tracked commit row with a
# newly-logged commit (avformat/avformat: also clear FFFormatContext packet
# queue when closing a muxer), re-sorted to the top of the "new entries"
# staging area (row 76) with a right/centre-aligned, date-formatted cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: drop the old last row (the hpeldsp_init_mmxext / 83b34691... row)
# This also shrinks the AutoFilter's data range down to A1:P75, matching the
# filter bookkeeping Excel leaves behind once that row is gone.
$ws.Rows.Item(76).Delete()

$ws.AutoFilterMode = $false
$ws.Range("A1:P75").AutoFilter() | Out-Null

# --- Step 2: write the newly logged commit into (the now-empty) row 76.
$ws.Range("A2").Copy()
$ws.Range("A76").PasteSpecial(-4122)

$ws.Range("A76").Value2 = 45674
$ws.Range("B76").Value2 = "c08d300481b8ebb846cd43a473988fdbc6793d1b"
$ws.Range("C76").Value2 = "avformat/avformat: also clear FFFormatContext packet queue when closing a muxer"
$ws.Range("D76").Value2 = "libavformat/avformat.c"
$ws.Range("E76").Value2 = "avformat_free_context"
$ws.Range("F76").Value2 = 966
$ws.Range("G76").Value2 = 50
$ws.Range("H76").Value2 = 1
$ws.Range("I76").Value2 = "Wrong"
$ws.Range("J76").Value2 = "Pass"
$ws.Range("M76").Value2 = 149
$ws.Range("N76").Value2 = 199
$ws.Range("O76").Value2 = 149
$ws.Range("P76").Value2 = 199

# Date cell keeps the yyyy-mm-dd format (copied above) but also picks up a
# right/vertically-centred alignment, matching the new combined style.
$ws.Range("A76").HorizontalAlignment = -4152
$ws.Range("A76").VerticalAlignment = -4108

# --- Step 3: the hidden _FilterDatabase name tracks the filtered range too.
$name = $wb.Names.Item(1)
$name.RefersTo = "=Sheet1!`$A`$1:`$P`$75"
